# Update LR-pairs data (Amelx-Cd63) with new TPM-derived values.
# Replaces the "MuSCs" sending-cluster block with an "ECs" block and
# refreshes every Ligand/Receptor/Edge expression metric (columns E:T)
# for rows 2-21 to match the new TPM computation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    @("ECs", "Amelx", "Cd63", "ECs", 1, 0.3333333333333333, 0.04979866666666666, 0.149396, 0.04800668898255549, 0.04800668898255549, 3, 1, 16.92841533333333, 50.785246, 0.0242040298661412, 0.0261568218313686, 0.8430125123795555, 7.587112611416, 0.001161955333908325, 0.00125570241043063),
    @("ECs", "Amelx", "Cd63", "FAPs", 1, 0.3333333333333333, 0.04979866666666666, 0.149396, 0.04800668898255549, 0.04800668898255549, 3, 1, 190.4708326666667, 571.412498, 0.2723327394629209, 0.2943046668003394, 9.485193505689777, 85.36674155120801, 0.01307379312316376, 0.01412859260519852),
    @("ECs", "Amelx", "Cd63", "Inflammatory-Mac", 1, 0.3333333333333333, 0.04979866666666666, 0.149396, 0.04800668898255549, 0.04800668898255549, 3, 1, 110.8005546666667, 332.401664, 0.1584212037293475, 0.17120269736802, 5.51771988832711, 49.659478994944, 0.007605277455676847, 0.00821887464552111),
    @("ECs", "Amelx", "Cd63", "MuSCs", 1, 0.3333333333333333, 0.04979866666666666, 0.149396, 0.04800668898255549, 0.04800668898255549, 2, 1, 156.6465685, 313.293137, 0.223971242892229, 0.1613608953572767, 7.800790249208666, 46.804741495252, 0.01075211779856363, 0.007746402317363465),
    @("ECs", "Amelx", "Cd63", "Resolving-Mac", 1, 0.3333333333333333, 0.04979866666666666, 0.149396, 0.04800668898255549, 0.04800668898255549, 3, 1, 224.5584563333333, 673.675369, 0.3210707840493613, 0.3469749186429952, 11.18271171412489, 100.644405427124, 0.01541354527124293, 0.01665711700404177),
    @("FAPs", "Amelx", "Cd63", "ECs", 3, 1, 0.809644, 2.428932, 0.7805094050963647, 0.7805094050963647, 3, 1, 16.92841533333333, 50.785246, 0.0242040298661412, 0.0261568218313686, 13.70598990414133, 123.353909137272, 0.01889147295175651, 0.02041564544681311),
    @("FAPs", "Amelx", "Cd63", "FAPs", 3, 1, 0.809644, 2.428932, 0.7805094050963647, 0.7805094050963647, 3, 1, 190.4708326666667, 571.412498, 0.2723327394629209, 0.2943046668003394, 154.2135668435707, 1387.922101592136, 0.2125582644664677, 0.2297075604014167),
    @("FAPs", "Amelx", "Cd63", "Inflammatory-Mac", 3, 1, 0.809644, 2.428932, 0.7805094050963647, 0.7805094050963647, 3, 1, 110.8005546666667, 332.401664, 0.1584212037293475, 0.17120269736802, 89.70900428253866, 807.3810385428479, 0.123649239477443, 0.1336253154736063),
    @("FAPs", "Amelx", "Cd63", "MuSCs", 3, 1, 0.809644, 2.428932, 0.7805094050963647, 0.7805094050963647, 2, 1, 156.6465685, 313.293137, 0.223971242892229, 0.1613608953572767, 126.827954306614, 760.967725839684, 0.174811661548507, 0.1259436964411248),
    @("FAPs", "Amelx", "Cd63", "Resolving-Mac", 3, 1, 0.809644, 2.428932, 0.7805094050963647, 0.7805094050963647, 3, 1, 224.5584563333333, 673.675369, 0.3210707840493613, 0.3469749186429952, 181.8124068195453, 1636.311661375908, 0.2505987666521904, 0.2708171873334038),
    @("Inflammatory-Mac", "Amelx", "Cd63", "ECs", 1, 0.3333333333333333, 0.036349, 0.109047, 0.0350410011879885, 0.0350410011879885, 3, 1, 16.92841533333333, 50.785246, 0.0242040298661412, 0.0261568218313686, 0.6153309689513333, 5.537978720562, 0.0008481334392935627, 0.0009165612248669906),
    @("Inflammatory-Mac", "Amelx", "Cd63", "FAPs", 1, 0.3333333333333333, 0.036349, 0.109047, 0.0350410011879885, 0.0350410011879885, 3, 1, 190.4708326666667, 571.412498, 0.2723327394629209, 0.2943046668003394, 6.923424296600666, 62.310818669406, 0.009542811847048372, 0.01031273017898125),
    @("Inflammatory-Mac", "Amelx", "Cd63", "Inflammatory-Mac", 1, 0.3333333333333333, 0.036349, 0.109047, 0.0350410011879885, 0.0350410011879885, 3, 1, 110.8005546666667, 332.401664, 0.1584212037293475, 0.17120269736802, 4.027489361578666, 36.247404254208, 0.005551237588082633, 0.005999113921859625),
    @("Inflammatory-Mac", "Amelx", "Cd63", "MuSCs", 1, 0.3333333333333333, 0.036349, 0.109047, 0.0350410011879885, 0.0350410011879885, 2, 1, 156.6465685, 313.293137, 0.223971242892229, 0.1613608953572767, 5.6939461184065, 34.163676710439, 0.007848176588261856, 0.005654247325909219),
    @("Inflammatory-Mac", "Amelx", "Cd63", "Resolving-Mac", 1, 0.3333333333333333, 0.036349, 0.109047, 0.0350410011879885, 0.0350410011879885, 3, 1, 224.5584563333333, 673.675369, 0.3210707840493613, 0.3469749186429952, 8.162475329260333, 73.46227796334301, 0.01125064172530207, 0.01215834853637141),
    @("Resolving-Mac", "Amelx", "Cd63", "ECs", 1, 0.3333333333333333, 0.141536, 0.424608, 0.1364429047330914, 0.1364429047330914, 3, 1, 16.92841533333333, 50.785246, 0.0242040298661412, 0.0261568218313686, 2.395980192618667, 21.563821733568, 0.003302468141182803, 0.003568912749257871),
    @("Resolving-Mac", "Amelx", "Cd63", "FAPs", 1, 0.3333333333333333, 0.141536, 0.424608, 0.1364429047330914, 0.1364429047330914, 3, 1, 190.4708326666667, 571.412498, 0.2723327394629209, 0.2943046668003394, 26.95847977230933, 242.626317950784, 0.03715787002624112, 0.04015578361474292),
    @("Resolving-Mac", "Amelx", "Cd63", "Inflammatory-Mac", 1, 0.3333333333333333, 0.141536, 0.424608, 0.1364429047330914, 0.1364429047330914, 3, 1, 110.8005546666667, 332.401664, 0.1584212037293475, 0.17120269736802, 15.68226730530133, 141.140405747712, 0.02161544920814503, 0.02335939332703304),
    @("Resolving-Mac", "Amelx", "Cd63", "MuSCs", 1, 0.3333333333333333, 0.141536, 0.424608, 0.1364429047330914, 0.1364429047330914, 2, 1, 156.6465685, 313.293137, 0.223971242892229, 0.1613608953572767, 22.171128719216, 133.026772315296, 0.03055928695689647, 0.02201654927287923),
    @("Resolving-Mac", "Amelx", "Cd63", "Resolving-Mac", 1, 0.3333333333333333, 0.141536, 0.424608, 0.1364429047330914, 0.1364429047330914, 3, 1, 224.5584563333333, 673.675369, 0.3210707840493613, 0.3469749186429952, 31.78310567559467, 286.047951080352, 0.04380783040062597, 0.04734226576917833)
)

for ($i = 0; $i -lt $newData.Length; $i++) {
    $rowValues = $newData[$i]
    $rowNum = $i + 2
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowValues[$j]
    }
}

Write-Output "Updated $($newData.Length) rows with new TPM data"
